$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing row 201 (current "Primera" /
# 44217 entry), shifting rows 201-211 down to 203-213. This mirrors the
# diff, where two brand-new "Especial" records (dated 44516) are spliced
# in at the top of the block and everything below cascades down by two
# rows.
$ws.Range("A201:A202").EntireRow.Insert()

# --- New row 201 ---------------------------------------------------
$ws.Range("A201").Value = 6
$ws.Range("B201").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C201").Value = "Metropolitana"
$ws.Range("D201").Value = 44516
$ws.Range("E201").Value = 13
$ws.Range("F201").Value = "Fruta"
$ws.Range("G201").Value = 100101
$ws.Range("H201").Value = "Berries"
$ws.Range("I201").Value = 100101001
$ws.Range("J201").Value = "Arándano (blue)"
$ws.Range("K201").Value = "Sin especificar"
$ws.Range("L201").Value = "Especial"
$ws.Range("M201").Value = 2500
$ws.Range("N201").Value = 6000
$ws.Range("O201").Value = 6000
$ws.Range("P201").Value = 6000
$ws.Range("Q201").Value = "$/bandeja 2 kilos"
$ws.Range("R201").Value = "Provincia de Curicó"
$ws.Range("S201").Value = 3000
$ws.Range("T201").Value = 2

# --- New row 202 ---------------------------------------------------
$ws.Range("A202").Value = 6
$ws.Range("B202").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C202").Value = "Metropolitana"
$ws.Range("D202").Value = 44516
$ws.Range("E202").Value = 13
$ws.Range("F202").Value = "Fruta"
$ws.Range("G202").Value = 100101
$ws.Range("H202").Value = "Berries"
$ws.Range("I202").Value = 100101001
$ws.Range("J202").Value = "Arándano (blue)"
$ws.Range("K202").Value = "Sin especificar"
$ws.Range("L202").Value = "Especial"
$ws.Range("M202").Value = 1250
$ws.Range("N202").Value = 5000
$ws.Range("O202").Value = 5000
$ws.Range("P202").Value = 5000
$ws.Range("Q202").Value = "$/bandeja 2 kilos"
$ws.Range("R202").Value = "Región de O'Higgins"
$ws.Range("S202").Value = 2500
$ws.Range("T202").Value = 2
